$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H33").Value = 1015.5
$ws.Range("I33").Value = 1015.5
$ws.Range("K33").Value = 1015.5
$ws.Range("M33").Value = -786.5
$ws.Range("H69").Value = 11665
$ws.Range("J69").Value = 11665
$ws.Range("L69").Value = 34995
$ws.Range("N69").Value = -36743
$ws.Range("H72").Value = 11665
$ws.Range("J72").Value = 11665
$ws.Range("L72").Value = 104985
$ws.Range("N72").Value = -113721
$ws.Range("H86").Value = 3512925.8
$ws.Range("I86").Value = 3700.4
$ws.Range("K86").Value = 3700.4
$ws.Range("M86").Value = -2577.4
$ws.Range("H87").Value = 81764.7
$ws.Range("I87").Value = 45000
$ws.Range("J87").Value = 89642.86
$ws.Range("K87").Value = 45000
$ws.Range("L87").Value = 89642.86
$ws.Range("M87").Value = -43752
$ws.Range("N87").Value = -92138.86
$ws.Range("H89").Value = 3512925.8
$ws.Range("I89").Value = 3700.4
$ws.Range("K89").Value = 18502
$ws.Range("M89").Value = -12886
$ws.Range("H90").Value = 81764.7
$ws.Range("I90").Value = 45000
$ws.Range("J90").Value = 89642.86
$ws.Range("K90").Value = 135000
$ws.Range("L90").Value = 268928.58
$ws.Range("M90").Value = -128760
$ws.Range("N90").Value = -281408.58
$ws.Range("H92").Value = 178
$ws.Range("I92").Value = 100.28571
$ws.Range("J92").Value = 450
$ws.Range("K92").Value = 100.28571
$ws.Range("L92").Value = 450
$ws.Range("M92").Value = 1147.71429
$ws.Range("N92").Value = -2946
$ws.Range("H94").Value = 3792.5
$ws.Range("I94").Value = 3723.3333
$ws.Range("J94").Value = 4000
$ws.Range("K94").Value = 3723.3333
$ws.Range("L94").Value = 4000
$ws.Range("M94").Value = -3272.3333
$ws.Range("N94").Value = -4902
$ws.Range("H96").Value = 627.6667
$ws.Range("I96").Value = 611.44446
$ws.Range("J96").Value = 676.3333
$ws.Range("K96").Value = 1834.33338
$ws.Range("L96").Value = 2028.9999
$ws.Range("M96").Value = -461.33338
$ws.Range("N96").Value = -4774.9999
$ws.Range("H98").Value = 2188
$ws.Range("I98").Value = 2044.4445
$ws.Range("J98").Value = 2762.2222
$ws.Range("K98").Value = 2044.4445
$ws.Range("L98").Value = 2762.2222
$ws.Range("M98").Value = -546.4445000000001
$ws.Range("N98").Value = -5758.2222
$ws.Range("H100").Value = 12306.846
$ws.Range("J100").Value = 14253.637
$ws.Range("L100").Value = 14253.637
$ws.Range("N100").Value = -15335.637
$ws.Range("H103").Value = 384.75
$ws.Range("I103").Value = 319.5
$ws.Range("K103").Value = 958.5
$ws.Range("M103").Value = -372.5
$ws.Range("H104").Value = 560
$ws.Range("I104").Value = 560
$ws.Range("K104").Value = 1680
$ws.Range("M104").Value = 67
$ws.Range("H106").Value = 3030.077
$ws.Range("I106").Value = 3030.077
$ws.Range("K106").Value = 3030.077
$ws.Range("M106").Value = -2399.077
$ws.Range("H112").Value = 3848.9614
$ws.Range("J112").Value = 3930.96
$ws.Range("L112").Value = 11792.88
$ws.Range("N112").Value = -14008.88
$ws.Range("H122").Value = 2188
$ws.Range("I122").Value = 2044.4445
$ws.Range("J122").Value = 2762.2222
$ws.Range("K122").Value = 6133.333500000001
$ws.Range("L122").Value = 8286.6666
$ws.Range("M122").Value = -3683.333500000001
$ws.Range("N122").Value = -13186.6666
$ws.Range("H129").Value = 1931.2307
$ws.Range("I129").Value = 544.5454999999999
$ws.Range("K129").Value = 1633.6365
$ws.Range("M129").Value = 3366.3635
$ws.Range("H132").Value = 1795.925
$ws.Range("I132").Value = 1851.5
$ws.Range("J132").Value = 740
$ws.Range("K132").Value = 5554.5
$ws.Range("L132").Value = 2220
$ws.Range("M132").Value = -3024.5
$ws.Range("N132").Value = -7280
$ws.Range("H135").Value = 572854.75
$ws.Range("I135").Value = 690925.9
$ws.Range("K135").Value = 6218333.100000001
$ws.Range("M135").Value = -6215798.100000001
$ws.Range("H137").Value = 4877.4375
$ws.Range("I137").Value = 4935.933
$ws.Range("K137").Value = 14807.799
$ws.Range("M137").Value = -12257.799
$ws.Range("H141").Value = 2457.05
$ws.Range("J141").Value = 1963
$ws.Range("L141").Value = 5889
$ws.Range("N141").Value = -16249
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 48799.953
$ws.Range("I2").Value = 53641.367
$ws.Range("K2").Value = 53641.367
$ws.Range("M2").Value = -53528.367
$ws.Range("H32").Value = 7023.1313
$ws.Range("I32").Value = 7091.757
$ws.Range("K32").Value = 7091.757
$ws.Range("M32").Value = -6804.757
$ws.Range("H61").Value = 2853.3076
$ws.Range("I61").Value = 2648
$ws.Range("K61").Value = 2648
$ws.Range("M61").Value = -2436
$ws.Range("H74").Value = 1647.9333
$ws.Range("I74").Value = 1647.9333
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1647.9333
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -773.9332999999999
$ws.Range("H77").Value = 1647.9333
$ws.Range("I77").Value = 1647.9333
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 8239.666499999999
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -3871.666499999999
$ws.Range("H102").Value = 1710
$ws.Range("I102").Value = 1780.1
$ws.Range("K102").Value = 1780.1
$ws.Range("M102").Value = -158.0999999999999
$ws.Range("H116").Value = 48799.953
$ws.Range("I116").Value = 53641.367
$ws.Range("K116").Value = 53641.367
$ws.Range("M116").Value = -51347.367
$ws.Range("H122").Value = 5004.472
$ws.Range("I122").Value = 4291.381
$ws.Range("K122").Value = 12874.143
$ws.Range("M122").Value = -10424.143
$ws.Range("H132").Value = 2043.3636
$ws.Range("I132").Value = 2075.2542
$ws.Range("K132").Value = 6225.7626
$ws.Range("M132").Value = -3695.7626
$ws.Range("H136").Value = 2853.3076
$ws.Range("I136").Value = 2648
$ws.Range("K136").Value = 7944
$ws.Range("M136").Value = -5394
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 48799.953
$ws.Range("I3").Value = 53641.367
$ws.Range("K3").Value = 53641.367
$ws.Range("M3").Value = -53527.367
$ws.Range("H105").Value = 84362
$ws.Range("I105").Value = 91848.45
$ws.Range("K105").Value = 91848.45
$ws.Range("M105").Value = -90101.45
$ws.Range("H107").Value = 358809.5
$ws.Range("I107").Value = 1292.6842
$ws.Range("J107").Value = 1113567.2
$ws.Range("K107").Value = 1292.6842
$ws.Range("L107").Value = 1113567.2
$ws.Range("M107").Value = 627.3158000000001
$ws.Range("N107").Value = -1117407.2
$ws.Range("H134").Value = 32456.771
$ws.Range("I134").Value = 4349.5356
$ws.Range("J134").Value = 144885.72
$ws.Range("K134").Value = 13048.6068
$ws.Range("L134").Value = 434657.16
$ws.Range("M134").Value = -10513.6068
$ws.Range("N134").Value = -439727.16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29709.395
$ws.Range("I31").Value = 2240.375
$ws.Range("J31").Value = 49686.863
$ws.Range("K31").Value = 2240.375
$ws.Range("L31").Value = 49686.863
$ws.Range("M31").Value = -1945.375
$ws.Range("N31").Value = -50276.863
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H34").Value = 29709.395
$ws.Range("I34").Value = 2240.375
$ws.Range("J34").Value = 49686.863
$ws.Range("K34").Value = 2240.375
$ws.Range("L34").Value = 49686.863
$ws.Range("M34").Value = -2038.375
$ws.Range("N34").Value = -50090.863
$ws.Range("H81").Value = 25000
$ws.Range("I81").Value = 25000
$ws.Range("K81").Value = 25000
$ws.Range("M81").Value = -24002
$ws.Range("H84").Value = 25000
$ws.Range("I84").Value = 25000
$ws.Range("K84").Value = 75000
$ws.Range("M84").Value = -70008
$ws.Range("H94").Value = 3099.875
$ws.Range("I94").Value = 3966.3333
$ws.Range("J94").Value = 2580
$ws.Range("K94").Value = 3966.3333
$ws.Range("L94").Value = 2580
$ws.Range("M94").Value = -3515.3333
$ws.Range("N94").Value = -3482
$ws.Range("H99").Value = 422170.1
$ws.Range("I99").Value = 4863
$ws.Range("J99").Value = 1006400
$ws.Range("K99").Value = 4863
$ws.Range("L99").Value = 1006400
$ws.Range("M99").Value = -3365
$ws.Range("N99").Value = -1009396
$ws.Range("H105").Value = 721.6667
$ws.Range("I105").Value = 721.6667
$ws.Range("K105").Value = 721.6667
$ws.Range("M105").Value = 1025.3333
$ws.Range("H107").Value = 324.5
$ws.Range("I107").Value = 327.22223
$ws.Range("K107").Value = 327.22223
$ws.Range("M107").Value = 1592.77777
$ws.Range("H126").Value = 422170.1
$ws.Range("I126").Value = 4863
$ws.Range("J126").Value = 1006400
$ws.Range("K126").Value = 14589
$ws.Range("L126").Value = 3019200
$ws.Range("M126").Value = -12119
$ws.Range("N126").Value = -3024140
$ws.Range("H134").Value = 252058.88
$ws.Range("I134").Value = 2215.4443
$ws.Range("J134").Value = 2500649.8
$ws.Range("K134").Value = 6646.3329
$ws.Range("L134").Value = 7501949.399999999
$ws.Range("M134").Value = -4111.3329
$ws.Range("N134").Value = -7507019.399999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 854.5625
$ws.Range("I14").Value = 854.5625
$ws.Range("K14").Value = 2563.6875
$ws.Range("M14").Value = -2390.6875
$ws.Range("H23").Value = 616.5
$ws.Range("J23").Value = 970.4286
$ws.Range("L23").Value = 2911.2858
$ws.Range("N23").Value = -3381.2858
$ws.Range("H34").Value = 83642.92999999999
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 83642.92999999999
$ws.Range("K34").Value = 0
$ws.Range("L34").ClearContents()
$ws.Range("M34").Value = 250928.79
$ws.Range("N34").Value = -251096.79
$ws.Range("H37").Value = 117727.86
$ws.Range("J37").Value = 117727.86
$ws.Range("L37").Value = 353183.58
$ws.Range("N37").Value = -353407.58
$ws.Range("H39").Value = 14186.875
$ws.Range("J39").Value = 27375
$ws.Range("L39").Value = 82125
$ws.Range("N39").Value = -82713
$ws.Range("H55").Value = 14646.363
$ws.Range("I55").Value = 5999.5
$ws.Range("J55").Value = 16567.889
$ws.Range("K55").Value = 17998.5
$ws.Range("L55").Value = 49703.667
$ws.Range("M55").Value = -17821.5
$ws.Range("N55").Value = -50057.667
$ws.Range("H56").Value = 7063.8
$ws.Range("I56").Value = 7063.8
$ws.Range("K56").Value = 7063.8
$ws.Range("M56").Value = -6533.8
$ws.Range("H103").Value = 4656.2856
$ws.Range("I103").Value = 36.666668
$ws.Range("J103").Value = 8121
$ws.Range("K103").Value = 110.000004
$ws.Range("L103").Value = 24363
$ws.Range("M103").Value = 768.999996
$ws.Range("N103").Value = -26121
$ws.Range("H122").Value = 46818.273
$ws.Range("J122").Value = 73067.28999999999
$ws.Range("L122").Value = 657605.61
$ws.Range("N122").Value = -662505.61
$ws.Range("H124").Value = 916
$ws.Range("I124").Value = 916
$ws.Range("K124").Value = 2748
$ws.Range("M124").Value = 2162
$ws.Range("H131").Value = 3641.4883
$ws.Range("J131").Value = 4299.5884
$ws.Range("L131").Value = 12898.7652
$ws.Range("N131").Value = -22978.7652
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 47247.5
$ws.Range("I40").Value = 44995
$ws.Range("K40").Value = 44995
$ws.Range("M40").Value = -44844
$ws.Range("H97").Value = 8960
$ws.Range("J97").Value = 11600
$ws.Range("L97").Value = 11600
$ws.Range("N97").Value = -12592
$ws.Range("H122").Value = 5727.8184
$ws.Range("J122").Value = 6142.857
$ws.Range("L122").Value = 18428.571
$ws.Range("N122").Value = -23328.571
$ws.Range("H126").Value = 3123.3872
$ws.Range("I126").Value = 2668.0557
$ws.Range("K126").Value = 8004.1671
$ws.Range("M126").Value = -5534.1671
$ws.Range("H132").Value = 41489.074
$ws.Range("I132").Value = 5512.278
$ws.Range("J132").Value = 113442.664
$ws.Range("K132").Value = 16536.834
$ws.Range("L132").Value = 340327.992
$ws.Range("M132").Value = -14006.834
$ws.Range("N132").Value = -345387.992
$ws.Range("H134").Value = 56583.168
$ws.Range("J134").Value = 56583.168
$ws.Range("L134").Value = 169749.504
$ws.Range("N134").Value = -174819.504
$ws.Range("H135").Value = 125075000
$ws.Range("J135").Value = 125075000
$ws.Range("L135").Value = 125075000
$ws.Range("N135").Value = -125085140
$ws.Range("H136").Value = 48862.23
$ws.Range("J136").Value = 48862.23
$ws.Range("L136").Value = 146586.69
$ws.Range("N136").Value = -151686.69
$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").ClearContents()
$ws.Range("N141").Value = 0
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1438089.2
$ws.Range("I7").Value = 17501.5
$ws.Range("K7").Value = 17501.5
$ws.Range("M7").Value = -17389.5
$ws.Range("H33").Value = 19999
$ws.Range("J33").Value = 19999
$ws.Range("L33").Value = 19999
$ws.Range("N33").Value = -20579
$ws.Range("H40").Value = 255501
$ws.Range("I40").Value = 503002
$ws.Range("K40").Value = 503002
$ws.Range("M40").Value = -502866
$ws.Range("H46").Value = 4041.625
$ws.Range("J46").Value = 6500
$ws.Range("L46").Value = 6500
$ws.Range("N46").Value = -6876
$ws.Range("H55").Value = 1083.4736
$ws.Range("I55").Value = 270.72726
$ws.Range("K55").Value = 270.72726
$ws.Range("M55").Value = -97.72726
$ws.Range("H61").Value = 6579.5
$ws.Range("I61").Value = 7266
$ws.Range("K61").Value = 7266
$ws.Range("M61").Value = -7064
$ws.Range("H93").Value = 2521.077
$ws.Range("I93").Value = 2097.7273
$ws.Range("K93").Value = 2097.7273
$ws.Range("M93").Value = -849.7273
$ws.Range("H113").Value = 6579.5
$ws.Range("I113").Value = 7266
$ws.Range("K113").Value = 7266
$ws.Range("M113").Value = -5096
$ws.Range("H122").Value = 1002078.06
$ws.Range("I122").Value = 557725.6
$ws.Range("K122").Value = 1673176.8
$ws.Range("M122").Value = -1670726.8
$ws.Range("H126").Value = 1438089.2
$ws.Range("I126").Value = 17501.5
$ws.Range("K126").Value = 52504.5
$ws.Range("M126").Value = -50034.5
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").ClearContents()
$ws.Range("N131").Value = 0
$ws.Range("H132").Value = 6387.222
$ws.Range("I132").Value = 5514.3
$ws.Range("J132").Value = 7478.375
$ws.Range("K132").Value = 16542.9
$ws.Range("L132").Value = 22435.125
$ws.Range("M132").Value = -14012.9
$ws.Range("N132").Value = -27495.125
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 7500
$ws.Range("J33").Value = 7500
$ws.Range("L33").Value = 7500
$ws.Range("N33").Value = -8000
$ws.Range("H36").Value = 7500
$ws.Range("J36").Value = 7500
$ws.Range("L36").Value = 7500
$ws.Range("N36").Value = -8000
$ws.Range("H40").Value = 89105.21000000001
$ws.Range("I40").Value = 93750
$ws.Range("J40").Value = 64333
$ws.Range("K40").Value = 93750
$ws.Range("L40").Value = 64333
$ws.Range("M40").Value = -93601
$ws.Range("N40").Value = -64631
$ws.Range("H62").Value = 82692
$ws.Range("I62").Value = 203999.2
$ws.Range("K62").Value = 203999.2
$ws.Range("M62").Value = -203375.2
$ws.Range("H65").Value = 82692
$ws.Range("I65").Value = 203999.2
$ws.Range("K65").Value = 1019996
$ws.Range("M65").Value = -1016876
$ws.Range("H74").Value = 13566.77
$ws.Range("I74").Value = 4997
$ws.Range("J74").Value = 15124.909
$ws.Range("K74").Value = 4997
$ws.Range("L74").Value = 15124.909
$ws.Range("M74").Value = -4061
$ws.Range("N74").Value = -16996.909
$ws.Range("H77").Value = 13566.77
$ws.Range("I77").Value = 4997
$ws.Range("J77").Value = 15124.909
$ws.Range("K77").Value = 14991
$ws.Range("L77").Value = 45374.727
$ws.Range("M77").Value = -10311
$ws.Range("N77").Value = -54734.727
$ws.Range("H113").Value = 100
$ws.Range("I113").Value = 100
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 300
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 1870
$ws.Range("H122").Value = 142861580
$ws.Range("I122").Value = 250002780
$ws.Range("J122").Value = 6667.6665
$ws.Range("K122").Value = 750008340
$ws.Range("L122").Value = 20002.9995
$ws.Range("M122").Value = -750005890
$ws.Range("N122").Value = -24902.9995
$ws.Range("H132").Value = 42836.883
$ws.Range("I132").Value = 3039.7222
$ws.Range("J132").Value = 132380.5
$ws.Range("K132").Value = 9119.1666
$ws.Range("L132").Value = 397141.5
$ws.Range("M132").Value = -6589.1666
$ws.Range("N132").Value = -402201.5
$ws.Range("H136").Value = 12335492
$ws.Range("I136").Value = 14929709
$ws.Range("K136").Value = 44789127
$ws.Range("M136").Value = -44786577
